$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1 (header): extend with two new columns P1, Q1 ---
# Copy formatting from O1 (bold + border, style index 1) onto the new header cells
$ws.Range("O1").Copy()
$ws.Range("P1:Q1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15

# --- Data rows 2-25: flip values in columns I, K, M, O and add P, Q ---
$ws.Range("I2:I25").Value = 2
$ws.Range("K2:K25").Value = 1
$ws.Range("M2:M25").Value = 2
$ws.Range("O2:O25").Value = 1
$ws.Range("P2:P25").Value = 2
$ws.Range("Q2:Q25").Value = 2
